$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the two anchor paragraphs that bound the block we need to
# rebuild: the "Separating the query processor..." heading through
# the "Using a VIEW..." heading (inclusive).
# ------------------------------------------------------------------
$sepIndex = -1
$viewIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($sepIndex -eq -1 -and $t -like "*Separating the query processor from the web UI*") {
        $sepIndex = $i
    }
    if ($t -like "*Using a VIEW, and not sending the result set as an array over the port*") {
        $viewIndex = $i
    }
}

if ($sepIndex -eq -1 -or $viewIndex -eq -1) {
    throw "Could not locate anchor paragraphs (sepIndex=$sepIndex viewIndex=$viewIndex)"
}

$startPara = $d.Paragraphs.Item($sepIndex)
$endPara = $d.Paragraphs.Item($viewIndex)
$fullRange = $d.Range($startPara.Range.Start, $endPara.Range.End)

# ------------------------------------------------------------------
# Replace that whole block with the new content:
#   - new "Not eliminating stop words when indexing" justification
#     (heading + body paragraph) inserted first
#   - the _GoBack bookmark now lives alone in its own paragraph
#   - the former content (Separating.../Honestly.../Even though...)
#     follows, with a lastRenderedPageBreak moved into the middle of
#     the "Honestly, this had to be done..." paragraph
#   - the "Using a VIEW..." heading no longer carries the page break
# ------------------------------------------------------------------
$xml = @"
<w:p>
  <w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr>
  <w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>Not eliminating stop words when indexing</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t>The elimination of stop-words had to be removed from the Indexer prior to the problem that would otherwise arise if the user types into the search bar a stop word, or if the phrase search of the user involves a stop word.</w:t></w:r>
</w:p>
<w:p>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr>
  <w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>Separating the query processor from the web UI</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">Honestly, this had to be done because the query processor is in Java, while the web UI is in HTML and CSS surely, but also in PHP. To solve this problem, a port was setup to allow for communication between the </w:t></w:r>
  <w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">web UI (PHP) and the query processor </w:t></w:r>
  <w:r><w:t>(Java), treating web UI (PHP) as the client and the query processor (Java) as the server.</w:t></w:r>
</w:p>
<w:p/>
<w:p>
  <w:r><w:t>Even though this solution was meant to solve the original problem of d</w:t></w:r>
  <w:r><w:t xml:space="preserve">ifferent languages, it made way </w:t></w:r>
  <w:r><w:t>the search engine more memory efficient, though unfortunately, not more time efficient.</w:t></w:r>
</w:p>
<w:p/>
<w:p>
  <w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr>
  <w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>Using a VIEW, and not sending the result set as an array over the port</w:t></w:r>
</w:p>
"@

$fullRange.InsertXML($xml)
